$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "testing" defined name so the sheet reference is single-quoted,
# e.g. 'testing'!$A$1:$H$12 instead of testing!$A$1:$H$12
$definedName = $wb.Names.Item(1)
$definedName.RefersTo = "='testing'!`$A`$1:`$H`$12"

# Refresh the date/datetime/time columns (D, E, G, H) with the full
# double-precision values instead of the previously rounded-off figures.
$ws.Range("G2").Value = 21916.000011574073
$ws.Range("H2").Value = 0.000011574074074074073
$ws.Range("D3").Value = 0.1
$ws.Range("G3").Value = 21916.000115740742
$ws.Range("H3").Value = 0.00003472222222222222
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 0.5
$ws.Range("G4").Value = 21916.00115740741
$ws.Range("H4").Value = 0.00010416666666666666
$ws.Range("D5").Value = 0.3
$ws.Range("E5").Value = 0.6666666666666666
$ws.Range("G5").Value = 21916.011574074077
$ws.Range("H5").Value = 0.0003125
$ws.Range("D6").Value = 0.4
$ws.Range("E6").Value = 0.75
$ws.Range("G6").Value = 21916.11574074074
$ws.Range("H6").Value = 0.0009375000000000001
$ws.Range("D7").Value = 0.5
$ws.Range("E7").Value = 0.8
$ws.Range("G7").Value = 21917.15740740741
$ws.Range("H7").Value = 0.0028125
$ws.Range("D8").Value = 0.6
$ws.Range("E8").Value = 0.8333333333333334
$ws.Range("G8").Value = 21927.574074074077
$ws.Range("H8").Value = 0.0084375
$ws.Range("D9").Value = 0.7
$ws.Range("E9").Value = 0.8571428571428571
$ws.Range("G9").Value = 22031.74074074074
$ws.Range("H9").Value = 0.0253125
$ws.Range("D10").Value = 0.8
$ws.Range("E10").Value = 0.875
$ws.Range("G10").Value = 23073.40740740741
$ws.Range("H10").Value = 0.07593749999999999
$ws.Range("D11").Value = 0.9
$ws.Range("E11").Value = 0.8888888888888888
$ws.Range("G11").Value = 33490.07407407407
$ws.Range("H11").Value = 0.2278125
$ws.Range("E12").Value = 0.9
$ws.Range("G12").Value = 137656.74074074073
$ws.Range("H12").Value = 0.6834375
